# Weekly refresh of the "Hortaliza, Mapocho Venta Directa de Santiago -
# Zapallo italiano" subset: rows 2-15 are re-aligned to a new weekly
# snapshot of the consolidated data (dates/volumes/prices/origins move to
# different rows; the descriptive columns A,B,C,E,F,G,H,I,R stay the same
# for every row since this sheet is a single market/variety slice).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44585; J=30; K=11000; L=11000; M=11000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=183; Q=60 },
    @{ Row=3;  D=44179; J=15; K=7000;  L=7000;  M=7000;  N='$/caja 60 unidades'; O='Provincia de Limarí';          P=117; Q=60 },
    @{ Row=4;  D=45001; J=40; K=10000; L=10000; M=10000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=167; Q=60 },
    @{ Row=5;  D=44186; J=15; K=7000;  L=7000;  M=7000;  N='$/caja 60 unidades'; O='Provincia de Limarí';          P=117; Q=60 },
    @{ Row=6;  D=44312; J=30; K=10000; L=10000; M=10000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=167; Q=60 },
    @{ Row=7;  D=44291; J=20; K=9000;  L=9000;  M=9000;  N='$/caja 60 unidades'; O='Provincia de Limarí';          P=150; Q=60 },
    @{ Row=8;  D=44405; J=45; K=9000;  L=9000;  M=9000;  N='$/caja 50 unidades'; O='Provincia de Quillota';        P=180; Q=50 },
    @{ Row=9;  D=44277; J=25; K=10000; L=10000; M=10000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=167; Q=60 },
    @{ Row=10; D=44284; J=35; K=10000; L=10000; M=10000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=167; Q=60 },
    @{ Row=11; D=44315; J=25; K=10000; L=10000; M=10000; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=167; Q=60 },
    @{ Row=12; D=44243; J=80; K=10000; L=11000; M=10375; N='$/caja 60 unidades'; O='Provincia de Quillota';        P=173; Q=60 },
    @{ Row=13; D=44333; J=25; K=10000; L=11000; M=10400; N='$/caja 60 unidades'; O='Provincia de Limarí';          P=173; Q=60 },
    @{ Row=14; D=45030; J=50; K=6000;  L=6000;  M=6000;  N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=120; Q=50 },
    @{ Row=15; D=44200; J=10; K=9000;  L=9000;  M=9000;  N='$/caja 60 unidades'; O='Provincia de Limarí';          P=150; Q=60 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
}
